$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 59, shifting existing rows 59-150 down to 62-153
$ws.Range("A59:A61").EntireRow.Insert()

# Row 59
$ws.Range('A59').Value = 10
$ws.Range('B59').Value = 'Vega Modelo de Temuco'
$ws.Range('C59').Value = 'La Araucanía'
$ws.Range('D59').Value = 44571
$ws.Range('E59').Value = 9
$ws.Range('F59').Value = 'Fruta'
$ws.Range('G59').Value = 100103
$ws.Range('H59').Value = 'Frutos de hueso (carozo)'
$ws.Range('I59').Value = 100103002
$ws.Range('J59').Value = 'Ciruela'
$ws.Range('K59').Value = 'Black Amber'
$ws.Range('L59').Value = 'Primera'
$ws.Range('M59').Value = 500
$ws.Range('N59').Value = 15000
$ws.Range('O59').Value = 16000
$ws.Range('P59').Value = 15600
$ws.Range('Q59').Value = '$/bandeja 18 kilos granel'
$ws.Range('R59').Value = 'Región de O''Higgins'
$ws.Range('S59').Value = 867
$ws.Range('T59').Value = 18

# Row 60
$ws.Range('A60').Value = 10
$ws.Range('B60').Value = 'Vega Modelo de Temuco'
$ws.Range('C60').Value = 'La Araucanía'
$ws.Range('D60').Value = 44571
$ws.Range('E60').Value = 9
$ws.Range('F60').Value = 'Fruta'
$ws.Range('G60').Value = 100103
$ws.Range('H60').Value = 'Frutos de hueso (carozo)'
$ws.Range('I60').Value = 100103002
$ws.Range('J60').Value = 'Ciruela'
$ws.Range('K60').Value = 'Black Amber'
$ws.Range('L60').Value = 'Primera'
$ws.Range('M60').Value = 8
$ws.Range('N60').Value = 350000
$ws.Range('O60').Value = 350000
$ws.Range('P60').Value = 350000
$ws.Range('Q60').Value = '$/bins (450 kilos)'
$ws.Range('R60').Value = 'Región de O''Higgins'
$ws.Range('S60').Value = 778
$ws.Range('T60').Value = 450

# Row 61
$ws.Range('A61').Value = 10
$ws.Range('B61').Value = 'Vega Modelo de Temuco'
$ws.Range('C61').Value = 'La Araucanía'
$ws.Range('D61').Value = 44571
$ws.Range('E61').Value = 9
$ws.Range('F61').Value = 'Fruta'
$ws.Range('G61').Value = 100103
$ws.Range('H61').Value = 'Frutos de hueso (carozo)'
$ws.Range('I61').Value = 100103002
$ws.Range('J61').Value = 'Ciruela'
$ws.Range('K61').Value = 'Black Amber'
$ws.Range('L61').Value = 'Segunda'
$ws.Range('M61').Value = 100
$ws.Range('N61').Value = 14000
$ws.Range('O61').Value = 14000
$ws.Range('P61').Value = 14000
$ws.Range('Q61').Value = '$/bandeja 18 kilos granel'
$ws.Range('R61').Value = 'Región de O''Higgins'
$ws.Range('S61').Value = 778
$ws.Range('T61').Value = 18
